# Weekly update: insert a new "Haba" price record for
# Terminal Hortofrutícola Agro Chillán and shift the existing
# historical rows (79-84) down by one (to 80-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 79; existing rows 79:84 shift to 80:85,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Cells.Item(79, 1).Value  = 7
$ws.Cells.Item(79, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(79, 3).Value  = "Ñuble"
$ws.Cells.Item(79, 4).Value2 = 45194
$ws.Cells.Item(79, 5).Value  = 16
$ws.Cells.Item(79, 6).Value  = 100112026
$ws.Cells.Item(79, 7).Value  = "Haba"
$ws.Cells.Item(79, 8).Value  = "Sin especificar"
$ws.Cells.Item(79, 9).Value  = "Primera"
$ws.Cells.Item(79, 10).Value = 60
$ws.Cells.Item(79, 11).Value = 14000
$ws.Cells.Item(79, 12).Value = 14000
$ws.Cells.Item(79, 13).Value = 14000
$ws.Cells.Item(79, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(79, 16).Value = 560
$ws.Cells.Item(79, 17).Value = 25
$ws.Cells.Item(79, 18).Value = "Hortaliza"
